# Refresh the crypto price/volume table (Sheet1, columns B-E, rows 2-51)
# with the latest scrape. Rows 16 and 17 swap contents (Polygon moves up
# to rank 14 / WrappedEther moves down to rank 15), so those two rows get
# their B/C/D/E cells rewritten in full rather than a simple value tweak.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "43.161.36"

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.314.34"
$ws.Range("E3").Value = "  +1.72%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  +0.03%  "

# --- Row 5: BNB ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.47"
$ws.Range("E5").Value = "  +1.40%  "

# --- Row 6: Solana ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.35"
$ws.Range("E6").Value = "  +5.33%  "

# --- Row 7: XRP ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.507"
$ws.Range("E7").Value = "  +2.87%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  -0.06%  "

# --- Row 9: Cardano ---
$ws.Range("E9").Value = "  +3.64%  "

# --- Row 10: Avalanche ---
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.56"
$ws.Range("E10").Value = "  +3.72%  "

# --- Row 11: Dogecoin ---
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("E11").Value = "  +0.99%  "

# --- Row 12: TRON ---
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.118"
$ws.Range("E12").Value = "  +4.19%  "

# --- Row 13: Chainlink ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.06"
$ws.Range("E13").Value = "  +14.67%  "

# --- Row 14: Polkadot ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("E14").Value = "  +3.50%  "

# --- Row 15: WrappedliquidstakedEther2.0 ---
$ws.Range("D15").Value = "2.691.54"
$ws.Range("E15").Value = "  +2.41%  "

# --- Row 16: Polygon (was WrappedEther) ---
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.819"
$ws.Range("E16").Value = "  +4.99%  "

# --- Row 17: WrappedEther (was Polygon) ---
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.269.25"
$ws.Range("E17").Value = "  +0.10%  "

# --- Row 18: WrappedBTC ---
$ws.Range("D18").Value = "43.112.08"
$ws.Range("E18").Value = "  +2.17%  "

# --- Row 19: InternetComputer(DFINITY) ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("E19").Value = "  +10.33%  "

# --- Row 20: ShibaInu ---
$ws.Range("D20").Value = "0.0₃0907"
$ws.Range("E20").Value = "  +1.97%  "

# --- Row 21: Uniswap ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.14"
$ws.Range("E21").Value = "  +2.30%  "

# --- Row 22: Litecoin ---
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.87"
$ws.Range("E22").Value = "  +1.81%  "

# --- Row 23: BitcoinCash ---
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.56"
$ws.Range("E23").Value = "  +1.89%  "

# --- Row 24: ImmutableX ---
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").Value = "  +12.84%  "

# --- Row 25: PancakeSwap ---
$ws.Range("E25").Value = "  +0.82%  "

# --- Row 26: Dai ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.31%  "

# --- Row 27: EthereumClassic ---
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.86"
$ws.Range("E27").Value = "  +4.12%  "

# --- Row 28: Monero ---
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.52"
$ws.Range("E28").Value = "  +0.40%  "

# --- Row 29: Toncoin ---
$ws.Range("E29").Value = "  -9.21%  "

# --- Row 30: InjectiveProtocol ---
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.19"
$ws.Range("E30").Value = "  +0.75%  "

# --- Row 31: Cosmos ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.20"
$ws.Range("E31").Value = "  +1.59%  "

# --- Row 32: FirstDigitalUSD ---
$ws.Range("E32").Value = "  +0.06%  "

# --- Row 33: Filecoin ---
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.05"
$ws.Range("E33").Value = "  +2.68%  "

# --- Row 34: RenderToken ---
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("E34").Value = "  +3.56%  "

# --- Row 35: WEMIXToken ---
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  +4.69%  "

# --- Row 36: Celestia ---
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.12"
$ws.Range("E36").Value = "  +5.83%  "

# --- Row 37: Hedera ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0693"
$ws.Range("E37").Value = "  +0.22%  "

# --- Row 38: Kaspa ---
$ws.Range("E38").Value = "  +3.99%  "

# --- Row 39: ARBITRUM ---
$ws.Range("E39").Value = "  +4.46%  "

# --- Row 40: LidoDAOToken ---
$ws.Range("E40").Value = "  +1.38%  "

# --- Row 41: Stellar ---
$ws.Range("E41").Value = "  +0.90%  "

# --- Row 42: ApeXProtocol ---
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("E42").Value = "  -1.79%  "

# --- Row 43: Maker ---
$ws.Range("D43").Value = "2.005.51"
$ws.Range("E43").Value = "  +2.24%  "

# --- Row 44: VeChain ---
$ws.Range("E44").Value = "  +3.63%  "

# --- Row 45: FraxShare ---
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.14"
$ws.Range("E45").Value = "  +5.55%  "

# --- Row 46: EnergySwap ---
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.62"
$ws.Range("E46").Value = "  +0.94%  "

# --- Row 47: NEARProtocol ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.87"
$ws.Range("E47").Value = "  +2.64%  "

# --- Row 48: MultiversX ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.81"
$ws.Range("E48").Value = "  +7.03%  "

# --- Row 49: RocketPoolETH ---
$ws.Range("D49").Value = "2.540.28"

# --- Row 50: Stacks ---
$ws.Range("E50").Value = "  +5.00%  "

# --- Row 51: THORChain ---
$ws.Range("E51").Value = "  +1.44%  "
